$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in Week 2 (D3) and the remaining-after-Janera (G3) values
$ws.Range("D3").Value = 1.5
$ws.Range("G3").Value = 0

# Row 5: fill in Week 2 (D5) and remaining (G5) values
$ws.Range("D5").Value = 3
$ws.Range("G5").Value = 0

# Row 13: fill in Week 2 (D13) and remaining (G13) values
$ws.Range("D13").Value = 8
$ws.Range("G13").Value = 0

# Row 18: new backlog item added by Janera
$ws.Range("A18").Value = "I want to be able to remove a meal from planned meals"
$ws.Range("B18").Value = "Implement functionality to remove a recipe from a specific day of the week for planned meals(desktop)"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = "Janera"
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 0

# Update the active cell selection to reflect where the author last left off
$ws.Range("C19").Select()

$wb.Save()
